# "Generate Report for handoff": the previous handoff attempt failed its
# transform, so the report is regenerated with a new source markdown file
# and a "Handoff transform failed" status, with the stale handoff-file /
# handoff-datetime / dependency-include info cleared out.

$wb = $excel.ActiveWorkbook

$newFile = "c3e54f12-3f2d-4217-be86-f2958ce1f92c.md"
$newStatus = "Handoff transform failed"
$epoch = "0001-01-01 00:00:00"

# --- Overview sheet ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A2").Value = $newFile
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus

# --- Per-language sheets (zh-cn, de-de) ---
$langSheets = @("zh-cn", "de-de")
foreach ($name in $langSheets) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A2").Value = $newFile
    $ws.Range("B2").Value = $newStatus

    # Handoff transform failed: no handoff file was produced this round,
    # so the latest-handoff-file cell (and its hyperlink formatting) is
    # cleared entirely rather than left blank.
    $ws.Range("C2").Clear()
    $ws.Range("D2").Value = $epoch

    # Dependency was ignored rather than included.
    $ws.Range("H2").Value = "Ignored"
}
